$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "22.072.17"
$ws.Range("E2").Value = "  -0.03%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.559.74"
$ws.Range("E3").Value = "  +0.53%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  -0.14%  "

# Row 5 - USDC
$ws.Range("D5").Value = "'0.9996"
$ws.Range("E5").Value = "  -0.25%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'291.92"

# Row 7 - XRP
$ws.Range("D7").Value = "'0.3972"
$ws.Range("E7").Value = "  +4.09%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3244"
$ws.Range("E8").Value = "  -0.89%  "

# Row 9 - OKB
$ws.Range("D9").Value = "'44.37"
$ws.Range("E9").Value = "  +1.80%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.07284"
$ws.Range("E10").Value = "  -0.83%  "

# Row 11 - Polygon
$ws.Range("E11").Value = "  -3.66%  "

# Row 12 - BinanceUSD
$ws.Range("D12").Value = "'0.9999"
$ws.Range("E12").Value = "  -0.15%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "'5.721"
$ws.Range("E13").Value = "  -0.81%  "

# Row 14 - Solana
$ws.Range("D14").Value = "'18.88"
$ws.Range("E14").Value = "  -5.77%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "'6.666"
$ws.Range("E15").Value = "  -1.06%  "

# Row 16 - was ShibaInu, now WrappedEther
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.564.02"
$ws.Range("E16").Value = "  +0.14%  "

# Row 17 - was WrappedEther, now ShibaInu
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.00001125"
$ws.Range("E17").Value = "  +4.44%  "

# Row 18 - TRON
$ws.Range("D18").Value = "'0.06597"
$ws.Range("E18").Value = "  -0.61%  "

# Row 19 - Litecoin
$ws.Range("D19").Value = "'83.94"
$ws.Range("E19").Value = "  -1.96%  "

# Row 20 - Dai
$ws.Range("D20").Value = "'0.9990"
$ws.Range("E20").Value = "  -0.33%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'6.295"
$ws.Range("E21").Value = "  -0.73%  "

# Row 22 - Avalanche
$ws.Range("D22").Value = "'15.62"
$ws.Range("E22").Value = "  -2.36%  "

# Row 23 - Cosmos
$ws.Range("E23").Value = "  -2.48%  "

# Row 24 - WrappedBTC
$ws.Range("D24").Value = "22.084.06"
$ws.Range("E24").Value = "  +0.07%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "'2.362"
$ws.Range("E25").Value = "  +2.61%  "

# Row 26 - LidoDAOToken
$ws.Range("D26").Value = "'2.431"
$ws.Range("E26").Value = "  -2.49%  "

# Row 27 - Monero
$ws.Range("D27").Value = "'148.46"
$ws.Range("E27").Value = "  -1.36%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'18.69"
$ws.Range("E28").Value = "  -2.29%  "

# Row 29 - HuobiToken
$ws.Range("D29").Value = "'4.877"
$ws.Range("E29").Value = "  -1.27%  "

# Row 30 - WrappedliquidstakedEther2.0
$ws.Range("D30").Value = "1.734.47"
$ws.Range("E30").Value = "  -0.03%  "

# Row 31 - BitcoinCash
$ws.Range("D31").Value = "'119.77"
$ws.Range("E31").Value = "  -1.51%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "'0.9906"
$ws.Range("E32").Value = "  -7.78%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "'5.955"
$ws.Range("E33").Value = "  +1.75%  "

# Row 34 - Stellar
$ws.Range("D34").Value = "'0.08339"
$ws.Range("E34").Value = "  +1.70%  "

# Row 35 - FraxShare
$ws.Range("D35").Value = "'9.185"
$ws.Range("E35").Value = "  -0.84%  "

# Row 36 - WEMIXTOKEN
$ws.Range("D36").Value = "'1.610"
$ws.Range("E36").Value = "  -15.32%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "'0.02284"
$ws.Range("E37").Value = "  -1.06%  "

# Row 38 - InternetComputer(DFINITY)
$ws.Range("D38").Value = "'5.162"
$ws.Range("E38").Value = "  -1.59%  "

# Row 39 - Hedera
$ws.Range("D39").Value = "'0.06036"
$ws.Range("E39").Value = "  -3.75%  "

# Row 40 - TrustWalletToken
$ws.Range("E40").Value = "  -1.39%  "

# Row 41 - Algorand
$ws.Range("D41").Value = "'0.2062"
$ws.Range("E41").Value = "  -3.87%  "

# Row 42 - Aptos
$ws.Range("E42").Value = "  -1.63%  "

# Row 43 - Frax
$ws.Range("D43").Value = "'0.9992"
$ws.Range("E43").Value = "  -0.29%  "

# Row 45 - was EnergySwap, now PancakeSwap
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").Value = "'3.769"
$ws.Range("E45").Value = "  +1.12%  "

# Row 46 - was PancakeSwap, now EnergySwap
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'13.06"
$ws.Range("E46").Value = "  -4.17%  "

# Row 47 - Decentraland
$ws.Range("D47").Value = "'0.5616"
$ws.Range("E47").Value = "  -3.42%  "

# Row 48 - Quant
$ws.Range("D48").Value = "'118.58"
$ws.Range("E48").Value = "  -2.54%  "

# Row 49 - NEARProtocol
$ws.Range("D49").Value = "'1.904"
$ws.Range("E49").Value = "  -3.05%  "

# Row 50 - EOS
$ws.Range("D50").Value = "'1.143"
$ws.Range("E50").Value = "  -2.33%  "

# Row 51 - Cronos
$ws.Range("D51").Value = "'0.06829"
$ws.Range("E51").Value = "  -2.53%  "
